$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 892.7143  # was 904.8333
$ws.Range("I28").Value = 909.75  # was 939.6667
$ws.Range("K28").Value = 909.75  # was 939.6667
$ws.Range("M28").Value = -424.75  # was -454.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 181.66667  # was 176.25
$ws.Range("I33").Value = 181.66667  # was 176.25
$ws.Range("K33").Value = 181.66667  # was 176.25
$ws.Range("M33").Value = 47.33332999999999  # was 52.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2363.6365  # was 2249.5833
$ws.Range("J43").Value = 2500  # was 2285
$ws.Range("L43").Value = 2500  # was 2285
$ws.Range("N43").Value = -2638  # was -2423

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1498  # was 0
$ws.Range("J127").Value = 1498  # was 0
$ws.Range("L127").Value = 4494  # was 0
$ws.Range("N127").Value = -14414  # newly added

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3517.3635  # was 3657.1667
$ws.Range("I137").Value = 3100.6667  # was 3250
$ws.Range("J137").Value = 3673.625  # was 3738.6
$ws.Range("K137").Value = 9302.000100000001  # was 9750
$ws.Range("L137").Value = 11020.875  # was 11215.8
$ws.Range("M137").Value = -6752.000100000001  # was -7200
$ws.Range("N137").Value = -16120.875  # was -16315.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 10665.733  # was 10665.866
$ws.Range("I138").Value = 999  # was 999.4
$ws.Range("K138").Value = 2997  # was 2998.2
$ws.Range("M138").Value = 2143  # was 2141.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 7999.5  # was 4999.75
$ws.Range("I141").Value = 7999.5  # was 4999.75
$ws.Range("K141").Value = 23998.5  # was 14999.25
$ws.Range("M141").Value = -18818.5  # was -9819.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9923.333000000001  # was 10874.737
$ws.Range("I32").Value = 9154.117  # was 10201.111
$ws.Range("K32").Value = 9154.117  # was 10201.111
$ws.Range("M32").Value = -8867.117  # was -9914.111000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3936.3333  # was 4054.625
$ws.Range("I61").Value = 3775.2856  # was 3906.1667
$ws.Range("K61").Value = 3775.2856  # was 3906.1667
$ws.Range("M61").Value = -3563.2856  # was -3694.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6236.8335  # was 6060.143
$ws.Range("J63").Value = 13245  # was 10496.667
$ws.Range("L63").Value = 13245  # was 10496.667
$ws.Range("N63").Value = -14617  # was -11868.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 6236.8335  # was 6060.143
$ws.Range("J66").Value = 13245  # was 10496.667
$ws.Range("L66").Value = 66225  # was 52483.335
$ws.Range("N66").Value = -73089  # was -59347.335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3936.3333  # was 4054.625
$ws.Range("I136").Value = 3775.2856  # was 3906.1667
$ws.Range("K136").Value = 11325.8568  # was 11718.5001
$ws.Range("M136").Value = -8775.856800000001  # was -9168.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 789.5  # was 4371.5
$ws.Range("J64").Value = 1234  # was 5713.6665
$ws.Range("L64").Value = 1234  # was 5713.6665
$ws.Range("N64").Value = -1684  # was -6163.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 789.5  # was 4371.5
$ws.Range("J67").Value = 1234  # was 5713.6665
$ws.Range("L67").Value = 1234  # was 5713.6665
$ws.Range("N67").Value = -2794  # was -7273.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 755  # was 689.5
$ws.Range("J80").Value = 1278  # was 1081.6666
$ws.Range("L80").Value = 1278  # was 1081.6666
$ws.Range("N80").Value = -3274  # was -3077.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 755  # was 689.5
$ws.Range("J83").Value = 1278  # was 1081.6666
$ws.Range("L83").Value = 6390  # was 5408.333000000001
$ws.Range("N83").Value = -16374  # was -15392.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 110000  # was 120000
$ws.Range("J112").Value = 110000  # was 120000
$ws.Range("L112").Value = 110000  # was 120000
$ws.Range("N112").Value = -112954  # was -122954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2068.8572  # was 2182.9473
$ws.Range("I31").Value = 1643.3334  # was 1744.6154
$ws.Range("K31").Value = 1643.3334  # was 1744.6154
$ws.Range("M31").Value = -1348.3334  # was -1449.6154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2068.8572  # was 2182.9473
$ws.Range("I34").Value = 1643.3334  # was 1744.6154
$ws.Range("K34").Value = 1643.3334  # was 1744.6154
$ws.Range("M34").Value = -1441.3334  # was -1542.6154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 19369.5  # was 19498.166
$ws.Range("I86").Value = 36000  # was 27748.5
$ws.Range("J86").Value = 2739  # was 2997.5
$ws.Range("K86").Value = 36000  # was 27748.5
$ws.Range("L86").Value = 2739  # was 2997.5
$ws.Range("M86").Value = -34877  # was -26625.5
$ws.Range("N86").Value = -4985  # was -5243.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 19369.5  # was 19498.166
$ws.Range("I89").Value = 36000  # was 27748.5
$ws.Range("J89").Value = 2739  # was 2997.5
$ws.Range("K89").Value = 180000  # was 138742.5
$ws.Range("L89").Value = 13695  # was 14987.5
$ws.Range("M89").Value = -174384  # was -133126.5
$ws.Range("N89").Value = -24927  # was -26219.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 60043.766  # was 78361.766
$ws.Range("I107").Value = 84437.086  # was 126400.5
$ws.Range("K107").Value = 84437.086  # was 126400.5
$ws.Range("M107").Value = -82517.086  # was -124480.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1117.1364  # was 1079.0416
$ws.Range("I5").Value = 897.55554  # was 854.36365
$ws.Range("K5").Value = 2692.66662  # was 2563.09095
$ws.Range("M5").Value = -2580.66662  # was -2451.09095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1150  # was 1360
$ws.Range("J68").Value = 1150  # was 1360
$ws.Range("L68").Value = 3450  # was 4080
$ws.Range("N68").Value = -5072  # was -5702

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1000  # was 1166.6666
$ws.Range("J69").Value = 0  # was 1500
$ws.Range("L69").Value = 0  # was 4500
$ws.Range("N69").Value = $null  # removed

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1000  # was 0
$ws.Range("I70").Value = 1000  # was 0
$ws.Range("K70").Value = 3000  # was 0
$ws.Range("M70").Value = -2685  # newly added

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1150  # was 1360
$ws.Range("J71").Value = 1150  # was 1360
$ws.Range("L71").Value = 10350  # was 12240
$ws.Range("N71").Value = -18462  # was -20352

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 1000  # was 1166.6666
$ws.Range("J72").Value = 0  # was 1500
$ws.Range("L72").Value = 0  # was 13500
$ws.Range("N72").Value = $null  # removed

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 1000  # was 0
$ws.Range("I73").Value = 1000  # was 0
$ws.Range("K73").Value = 3000  # was 0
$ws.Range("M73").Value = -1908  # newly added

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 7216.6665  # was 6812.5
$ws.Range("J74").Value = 7216.6665  # was 6812.5
$ws.Range("L74").Value = 21649.9995  # was 20437.5
$ws.Range("N74").Value = -23771.9995  # was -22559.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 0  # was 156.5
$ws.Range("I75").Value = 0  # was 156.5
$ws.Range("K75").Value = 0  # was 469.5
$ws.Range("M75").Value = $null  # removed

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 7216.6665  # was 6812.5
$ws.Range("J77").Value = 7216.6665  # was 6812.5
$ws.Range("L77").Value = 64949.9985  # was 61312.5
$ws.Range("N77").Value = -75557.9985  # was -71920.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 0  # was 156.5
$ws.Range("I78").Value = 0  # was 156.5
$ws.Range("K78").Value = 0  # was 1408.5
$ws.Range("M78").Value = $null  # removed

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 400  # was 450
$ws.Range("I86").Value = 400  # was 450
$ws.Range("K86").Value = 1200  # was 1350
$ws.Range("M86").Value = -14  # was -164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 5599.6665  # was 8249.5
$ws.Range("I87").Value = 8149.5  # was 15999
$ws.Range("K87").Value = 24448.5  # was 47997
$ws.Range("M87").Value = -23200.5  # was -46749

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 400  # was 450
$ws.Range("I89").Value = 400  # was 450
$ws.Range("K89").Value = 3600  # was 4050
$ws.Range("M89").Value = 2328  # was 1878

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 5599.6665  # was 8249.5
$ws.Range("I90").Value = 8149.5  # was 15999
$ws.Range("K90").Value = 73345.5  # was 143991
$ws.Range("M90").Value = -67105.5  # was -137751

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 549.875  # was 585.7143
$ws.Range("I107").Value = 399.83334  # was 420
$ws.Range("K107").Value = 1199.50002  # was 1260
$ws.Range("M107").Value = 720.4999800000001  # was 660

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1117.1364  # was 1079.0416
$ws.Range("I135").Value = 897.55554  # was 854.36365
$ws.Range("K135").Value = 8077.99986  # was 7689.27285
$ws.Range("M135").Value = -5542.99986  # was -5154.27285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 7500  # was 7666
$ws.Range("I97").Value = 4000  # was 5999
$ws.Range("K97").Value = 4000  # was 5999
$ws.Range("M97").Value = -3504  # was -5503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 773.3333  # was 774.5
$ws.Range("I102").Value = 773.3333  # was 774.5
$ws.Range("K102").Value = 773.3333  # was 774.5
$ws.Range("M102").Value = 848.6667  # was 847.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2401.4285  # was 2320.5
$ws.Range("J113").Value = 2000  # was 2065.8333
$ws.Range("L113").Value = 2000  # was 2065.8333
$ws.Range("N113").Value = -6340  # was -6405.8333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1984.3846  # was 2208.6667
$ws.Range("I122").Value = 2299.6667  # was 2489.3333
$ws.Range("J122").Value = 1275  # was 1366.6666
$ws.Range("K122").Value = 6899.000100000001  # was 7467.999899999999
$ws.Range("L122").Value = 3825  # was 4099.9998
$ws.Range("M122").Value = -4449.000100000001  # was -5017.999899999999
$ws.Range("N122").Value = -8725  # was -8999.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 552.5714  # was 566.2857
$ws.Range("I55").Value = 631.0909  # was 648.5454999999999
$ws.Range("K55").Value = 631.0909  # was 648.5454999999999
$ws.Range("M55").Value = -458.0909  # was -475.5454999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 0  # was 19000
$ws.Range("I94").Value = 0  # was 19000
$ws.Range("K94").Value = 0  # was 19000
$ws.Range("M94").Value = $null  # removed

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2197.8  # was 3000
$ws.Range("I122").Value = 997.25  # was 1000
$ws.Range("K122").Value = 2991.75  # was 3000
$ws.Range("M122").Value = -541.75  # was -550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4856.4287  # was 4999.5
$ws.Range("I132").Value = 3999  # was 3999.5
$ws.Range("K132").Value = 11997  # was 11998.5
$ws.Range("M132").Value = -9467  # was -9468.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 606.3333  # was 662.6
$ws.Range("I107").Value = 635.1111  # was 605.1429000000001
$ws.Range("J107").Value = 563.1667  # was 796.6667
$ws.Range("K107").Value = 1905.3333  # was 1815.4287
$ws.Range("L107").Value = 1689.5001  # was 2390.0001
$ws.Range("M107").Value = 14.66670000000022  # was 104.5712999999998
$ws.Range("N107").Value = -5529.5001  # was -6230.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 821.06665  # was 822.8
$ws.Range("I113").Value = 742.8333  # was 708.1429000000001
$ws.Range("J113").Value = 873.2222  # was 923.125
$ws.Range("K113").Value = 2228.4999  # was 2124.4287
$ws.Range("L113").Value = 2619.6666  # was 2769.375
$ws.Range("M113").Value = -58.4998999999998  # was 45.57129999999961
$ws.Range("N113").Value = -6959.6666  # was -7109.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4666  # was 2747.8333
$ws.Range("I132").Value = 4000  # was 1622.25
$ws.Range("K132").Value = 12000  # was 4866.75
$ws.Range("M132").Value = -9470  # was -2336.75
